$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '20220603181019556_boot'
$ws.Range("D2").Value = 'ServerError'
$ws.Range("H2").Value = '<Response [500]>'
$ws.Range("C3").Value = '20220603181019556_boot'
$ws.Range("D3").Value = 'ServerError'
$ws.Range("G3").Value = '{''connectorId'': ''01'', ''errorCode'': None, ''info'': [{''reason'': None, ''cpv'': 100, ''rv'': 11}], ''status'': ''Available'', ''timestamp'': ''2022-06-03T18:10:20Z'', ''vendorErrorCode'': '''', ''vendorId'': ''LGE''}'
$ws.Range("H3").Value = '<Response [500]>'
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '1111222233334444'
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = '20220603181021746_card'
$ws.Range("G4").Value = '{''idTag'': ''1111222233334444''}'
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '1111222233334444'
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = '20220603181021746_card'
$ws.Range("D5").Value = 'ServerError'
$ws.Range("G5").Value = '{''connectorId'': ''01'', ''errorCode'': None, ''info'': [{''reason'': None, ''cpv'': 100, ''rv'': 11}], ''status'': ''Preparing'', ''timestamp'': ''2022-06-03T18:10:22Z'', ''vendorErrorCode'': '''', ''vendorId'': ''LGE''}'
$ws.Range("H5").Value = '<Response [500]>'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '1111222233334444'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = '20220603181021746_card'
$ws.Range("G6").Value = '{''venderId'': ''LG'', ''messageId'': ''Tariff'', ''data'': {''connectorId'': ''01'', ''idTag'': ''1111222233334444'', ''timestamp'': ''2022-06-03T18:10:23Z''}}'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '1111222233334444'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '20220603181021746_card'
$ws.Range("G7").Value = '{''idTag'': ''5555222233334444'', ''connectorId'': ''01'', ''meterStart'': None, ''timestamp'': ''2022-06-03T18:10:24Z''}'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '1111222233334444'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '20220603181021746_card'
$ws.Range("D8").Value = 'ServerError'
$ws.Range("G8").Value = '{''connectorId'': ''01'', ''errorCode'': None, ''info'': [{''reason'': None, ''cpv'': 100, ''rv'': 11}], ''status'': ''Charging'', ''timestamp'': ''2022-06-03T18:10:26Z'', ''vendorErrorCode'': '''', ''vendorId'': ''LGE''}'
$ws.Range("H8").Value = '<Response [500]>'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '1111222233334444'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = '20220603181021746_card'
$ws.Range("D9").Value = 'ServerError'
$ws.Range("G9").Value = '{''connectorId'': ''01'', ''transactionId'': None, ''meterValue'': [{''timestamp'': ''2022-06-03T18:10:27Z'', ''sampledValue'': [{''measurand'': ''Current.Import'', ''phase'': ''L1'', ''unit'': ''A'', ''value'': ''23.4''}, {''measurand'': ''Voltage'', ''phase'': ''L1'', ''unit'': ''V'', ''value'': ''220.7''}, {''measurand'': ''Energy.Active.Import.Register'', ''unit'': ''Wh'', ''value'': ''999.8''}, {''measurand'': ''SoC'', ''unit'': ''%'', ''value'': ''10''}, {''measurand'': ''Power.Active.Import'', ''unit'': ''W'', ''value'': ''0.7''}]}]}'
$ws.Range("H9").Value = '<Response [500]>'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '1111222233334444'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = '20220603181021746_card'
$ws.Range("G10").Value = '{''idTag'': ''3333222233334444'', ''meterStop'': 0.729759527533703, ''reason'': ''Finished'', ''timestamp'': ''2022-06-03T18:10:28Z'', ''transactionId'': None, ''transactionData'': [{''timestamp'': ''2022-06-03T18:10:28Z'', ''sampledValue'': [{''measurand'': ''01'', ''phase'': ''01'', ''unit'': ''01'', ''value'': ''01''}, {''measurand'': ''01'', ''phase'': ''01'', ''unit'': ''01'', ''value'': ''01''}]}]}'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '1111222233334444'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = '20220603181021746_card'
$ws.Range("D11").Value = 'ServerError'
$ws.Range("G11").Value = '{''connectorId'': ''01'', ''errorCode'': None, ''info'': [{''reason'': None, ''cpv'': 100, ''rv'': 11}], ''status'': ''Finishing'', ''timestamp'': ''2022-06-03T18:10:29Z'', ''vendorErrorCode'': '''', ''vendorId'': ''LGE''}'
$ws.Range("H11").Value = '<Response [500]>'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '1111222233334444'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = '20220603181021746_card'
$ws.Range("D12").Value = 'ServerError'
$ws.Range("G12").Value = '{''connectorId'': ''01'', ''errorCode'': None, ''info'': [{''reason'': None, ''cpv'': 100, ''rv'': 11}], ''status'': ''Available'', ''timestamp'': ''2022-06-03T18:10:30Z'', ''vendorErrorCode'': '''', ''vendorId'': ''LGE''}'
$ws.Range("H12").Value = '<Response [500]>'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '1111222233334444'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = '20220603181034632_heartbeat'
